# Applies the "Updated cryptos list" data refresh: new Price (column D) and
# Volume(1h) (column E) figures for (almost) every coin row, plus a couple of
# rows (POPCAT/InternetComputer, Bittensor/FirstDigitalUSD) that swapped rank
# order and therefore swapped Coin/Link/Price/Volume content between rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style of an untouched text data cell (no explicit "s" attribute). We reapply
# this to every Price cell we touch so that forcing it to remain text (via a
# leading apostrophe, see below) does not leave a stray quote-prefix/
# text-number-format style attached to the cell.
$refStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '''68.497.46'
$ws.Range("D2").Style = $refStyle
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("D3").Value = '''2.459.07'
$ws.Range("D3").Style = $refStyle
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("D5").Value = '''559.31'
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = '  -2.60%  '
$ws.Range("D6").Value = '''164.23'
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = '  -1.35%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.504'
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = '  -1.62%  '
$ws.Range("D9").Value = '''2.458.76'
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = '  -1.76%  '
$ws.Range("D10").Value = '''0.150'
$ws.Range("D10").Style = $refStyle
$ws.Range("E10").Value = '  -6.01%  '
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("D12").Value = '''0.336'
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = '  -6.14%  '
$ws.Range("D13").Value = '''4.82'
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").Value = '''2.911.24'
$ws.Range("D14").Style = $refStyle
$ws.Range("E14").Value = '  -1.61%  '
$ws.Range("D15").Value = '''68.437.87'
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("E16").Value = '  -3.81%  '
$ws.Range("D17").Value = '''23.45'
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = '  -5.09%  '
$ws.Range("D18").Value = '''2.473.38'
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = '  -1.15%  '
$ws.Range("E19").Value = '  -1.77%  '
$ws.Range("D20").Value = '''344.41'
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = '  -1.30%  '
$ws.Range("D21").Value = '''7.17'
$ws.Range("D21").Style = $refStyle
$ws.Range("E21").Value = '  -4.48%  '
$ws.Range("D22").Value = '''3.77'
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = '  -3.34%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("D24").Value = '''1.87'
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = '  -3.67%  '
$ws.Range("D25").Value = '''67.78'
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = '  -4.28%  '
$ws.Range("D26").Value = '''1.11'
$ws.Range("D26").Style = $refStyle
$ws.Range("E26").Value = '  +11.67%  '
$ws.Range("D27").Value = '''3.71'
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = '  -5.50%  '
$ws.Range("D29").Value = '''8.13'
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = '  -6.89%  '
$ws.Range("D30").Value = '''0.0₃0832'
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = '  -6.47%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''7.21'
$ws.Range("D31").Style = $refStyle
$ws.Range("E31").Value = '  -8.36%  '
$ws.Range("B32").Value = 'POPCAT'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range("D32").Value = '''3.41'
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = '  +131.91%  '
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = '''1.00'
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Value = '''433.04'
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = '  -5.45%  '
$ws.Range("D35").Value = '''1.16'
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = '  -3.33%  '
$ws.Range("D36").Value = '''1.67'
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = '  -3.60%  '
$ws.Range("D37").Value = '''157.48'
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("D39").Value = '''0.110'
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = '  -4.62%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").Value = '''17.90'
$ws.Range("D41").Style = $refStyle
$ws.Range("E41").Value = '  -2.45%  '
$ws.Range("E42").Value = '  -3.54%  '
$ws.Range("D43").Value = '''4.46'
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = '  -4.84%  '
$ws.Range("E44").Value = '  -5.29%  '
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").Value = '''2.08'
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = '  -5.37%  '
$ws.Range("D47").Value = '''134.86'
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = '  -4.48%  '
$ws.Range("D48").Value = '''3.35'
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = '  -3.77%  '
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("D50").Value = '''0.483'
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = '  -6.86%  '
$ws.Range("D51").Value = '''0.564'
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = '  -2.56%  '
